$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "FISHING PIERS dist" table (columns P:R), mirrors the existing
# PRIVATE SCHOOLS dist / TRAFFIC CAMERAS dist tables in columns G:I / K:N
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "FISHING PIERS dist"
$ws.Range("Q1").Value = "classifier"
$ws.Range("R1").Value = "error"

$ws.Range("P2").Value = 8
$ws.Range("Q2").Value = "random forest"
$ws.Range("R2").Value = 58.41

$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = "random forest"
$ws.Range("R3").Value = 58.45

$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = "random forest"
$ws.Range("R4").Value = 58.31

$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = "random forest"
$ws.Range("R5").Value = 58.35

$ws.Range("P6").Value = 0.5
$ws.Range("Q6").Value = "random forest"
$ws.Range("R6").Value = 58.51

$ws.Range("P7").Value = 0.25
$ws.Range("Q7").Value = "random forest"
$ws.Range("R7").Value = 58.53

# ---------------------------------------------------------------------------
# New "MOVIE THEATERS" table (columns T:V), same shape as FISHING PIERS dist
# ---------------------------------------------------------------------------
$ws.Range("T1").Value = "MOVIE THEATERS"
$ws.Range("U1").Value = "classifier"
$ws.Range("V1").Value = "error"

$ws.Range("T2").Value = 8
$ws.Range("U2").Value = "random forest"
$ws.Range("V2").Value = 58.53

$ws.Range("T3").Value = 4
$ws.Range("U3").Value = "random forest"
$ws.Range("V3").Value = 58.43

$ws.Range("T4").Value = 2
$ws.Range("U4").Value = "random forest"
$ws.Range("V4").Value = 58.34

$ws.Range("T5").Value = 1
$ws.Range("U5").Value = "random forest"
$ws.Range("V5").Value = 58.38

$ws.Range("T6").Value = 0.5
$ws.Range("U6").Value = "random forest"
$ws.Range("V6").Value = 58.39

$ws.Range("T7").Value = 0.25
$ws.Range("U7").Value = "random forest"
$ws.Range("V7").Value = 58.59

# Header row wraps text like the other header cells in row 1 (style "s=2")
$ws.Range("P1:R1,T1:V1").WrapText = $true

# ---------------------------------------------------------------------------
# New "Fishing piers = 2" mini comparison table (columns P:S, row 12 header)
# ---------------------------------------------------------------------------
$ws.Range("P12").Value = "Fishing piers = 2"
$ws.Range("Q12").Value = "random forest"
$ws.Range("R12").Value = "SVC (linear)"
$ws.Range("S12").Value = "SVC(rbf) 100 folds"
$ws.Range("P12,S12").WrapText = $true

$ws.Range("Q13").Value = 58.31
$ws.Range("R13").Value = 77.48
$ws.Range("S13").Value = 68.55
$ws.Range("S14").Value = 67.89
$ws.Range("S15").Value = 67.02
$ws.Range("S16").Value = 67.34

# Bold emphasis, matching the existing "best value" styling in this sheet
$ws.Range("Q13").Font.Bold = $true

# ---------------------------------------------------------------------------
# Bold emphasis on existing cells that got highlighted in the same edit
# ---------------------------------------------------------------------------
$ws.Range("A4").Font.Bold = $true
$ws.Range("D4").Font.Bold = $true
$ws.Range("B25").Font.Bold = $true
$ws.Range("C25").Font.Bold = $true

# ---------------------------------------------------------------------------
# New trailing rows: a small "WITH ALL FEATURES" note block
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = "error"

$ws.Range("A29").Value = 15
$ws.Range("B29").Value = 80
$ws.Range("C29").Value = 58.2
$ws.Range("C29").Font.Bold = $true
$ws.Range("D29").Value = "WITH ALL FEATURES"

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 16.92
$ws.Columns.Item(17).ColumnWidth = 14.17
$ws.Columns.Item(21).ColumnWidth = 14.17

# ---------------------------------------------------------------------------
# Selection matches the author's final cursor position
# ---------------------------------------------------------------------------
$null = $ws.Range("D29").Select()

Write-Host "edit applied"
